$p = $ppt.ActivePresentation

# Slide 5 contains a single table (the "Type of document / Definition / Why it
# is important" table). The author switched the table to a different built-in
# PowerPoint table style ("Medium Style 2 - Accent 1",
# {21BE59EB-6F66-4F38-8893-DA6D40029C92}).
$s = $p.Slides.Item(5)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{21BE59EB-6F66-4F38-8893-DA6D40029C92}")
    }
}
